$wb = $excel.ActiveWorkbook

# --- Update the "invalid_login" sheet (sheet3 / rId3) with new test data ---
$ws = $wb.Worksheets.Item("invalid_login")

# Write the brand-new unique strings in the same order they first appear in
# the final shared-strings table: testrecep, test_password, deneme.
$ws.Range("A3").Value = "testrecep"
$ws.Range("B4").Value = "test_password"
$ws.Range("B2").Value = "deneme"

# Fill in the remaining (already-known) shared strings for the new layout.
$ws.Range("A2").Value = "recepodemis"
$ws.Range("B3").Value = "secret_sauce"
$ws.Range("A4").Value = "standard_user"

# Move the selection/active cell on this sheet to D6.
$ws.Activate() | Out-Null
$ws.Range("D6").Select() | Out-Null

# --- Make "locked_out" (sheet4 / rId4) the active/selected tab ---
$ws2 = $wb.Worksheets.Item("locked_out")
$ws2.Activate() | Out-Null
$ws2.Range("D7").Select() | Out-Null
